# Log.docx update:
#  - add a "gramStart"/"gramEnd" proofErr pair around the  ”  ->”  ... “.”  .
#    sentence in the first (existing) log entry
#  - split the paragraph after "... forma correcta." and append a brand new
#    log entry (14/02/2014) with its own spell/gram proofErr markers,
#    moving the _GoBack bookmark to the end of that new entry

$d = $word.ActiveDocument

# First paragraph (the single existing log entry) - grab it dynamically so
# we are not relying on a hard-coded character offset.
$p1 = $d.Paragraphs(1)
$p1Range = $p1.Range
$textEnd = $p1Range.End - 1   # End is just past the paragraph mark; back up one.

$replaceRange = $d.Range($p1Range.Start, $textEnd)

$xml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="003D1985" w:rsidRDefault="006A624A"><w:r><w:t>13/02/2014 2</w:t></w:r><w:r w:rsidR="00D00665"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00D00665"><w:t>hr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="009919EF"><w:t xml:space="preserve">. Modificadas las </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="009919EF"><w:t>Invalid</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="009919EF"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="009919EF"><w:t>Expression</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="009919EF"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="009919EF"><w:t>Exception</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="009919EF"><w:t xml:space="preserve"> y agregada aclaración que se debe usar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="009919EF"><w:t>Typedef</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="009919EF"><w:t xml:space="preserve"> en las pruebas.</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>”</w:t></w:r><w:r w:rsidR="00A7169F"><w:t xml:space="preserve"> -&gt;</w:t></w:r><w:r><w:t>”</w:t></w:r><w:r w:rsidR="00A7169F"><w:t xml:space="preserve"> será tratado como </w:t></w:r><w:r><w:t>“</w:t></w:r><w:r w:rsidR="00A7169F"><w:t>.</w:t></w:r><w:r><w:t>”</w:t></w:r><w:r w:rsidR="00A7169F"><w:t xml:space="preserve"> .</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> Ahora puede leer varios atributos declarados en la misma línea de forma correcta.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">14/02/2014 3hr. Ya guarda los atributos </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>globales</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>,pero</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> todavía no los incluye en las heurísticas ni los muestra. No corre todavía con </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>Aeropuerto.c</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> .</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> Ahora lee números de </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>mas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> de 1 cifra. Ahora puede leer caracteres encerrados entre apostrofes. Ahora soporta funciones que devuelven estructuras.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$replaceRange.InsertXML($xml)
